$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 13150
$ws.Range("J21").Value = 18750
$ws.Range("L21").Value = 18750
$ws.Range("N21").Value = -19686
$ws.Range("H23").Value = 13150
$ws.Range("J23").Value = 18750
$ws.Range("L23").Value = 18750
$ws.Range("N23").Value = -19218
$ws.Range("H28").Value = 654249.9399999999
$ws.Range("I28").Value = 1010350.06
$ws.Range("J28").Value = 1399.6666
$ws.Range("K28").Value = 1010350.06
$ws.Range("L28").Value = 1399.6666
$ws.Range("M28").Value = -1009865.06
$ws.Range("N28").Value = -2369.6666
$ws.Range("H34").Value = 3972
$ws.Range("I34").Value = 3972
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3972
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3769
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 3972
$ws.Range("I36").Value = 3972
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 3972
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -3257
$ws.Range("N36").ClearContents()
$ws.Range("H62").Value = 5860271.5
$ws.Range("I62").Value = 10111174
$ws.Range("K62").Value = 10111174
$ws.Range("M62").Value = -10110550
$ws.Range("H65").Value = 5860271.5
$ws.Range("I65").Value = 10111174
$ws.Range("K65").Value = 50555870
$ws.Range("M65").Value = -50552750
$ws.Range("H107").Value = 654187.8
$ws.Range("I107").Value = 926557.5600000001
$ws.Range("J107").Value = 500.4
$ws.Range("K107").Value = 926557.5600000001
$ws.Range("L107").Value = 500.4
$ws.Range("M107").Value = -924637.5600000001
$ws.Range("N107").Value = -4340.4
$ws.Range("H111").Value = 1396.8572
$ws.Range("I111").Value = 957.7
$ws.Range("J111").Value = 2494.75
$ws.Range("K111").Value = 2873.1
$ws.Range("L111").Value = 7484.25
$ws.Range("M111").Value = 193.8999999999996
$ws.Range("N111").Value = -13618.25
$ws.Range("H116").Value = 11533331
$ws.Range("I116").Value = 17298812
$ws.Range("J116").Value = 2369.25
$ws.Range("K116").Value = 17298812
$ws.Range("L116").Value = 2369.25
$ws.Range("M116").Value = -17295370
$ws.Range("N116").Value = -9253.25
$ws.Range("H118").Value = 422
$ws.Range("I118").Value = 422
$ws.Range("K118").Value = 1266
$ws.Range("M118").Value = 391

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 80377.30499999999
$ws.Range("I2").Value = 147958.86
$ws.Range("J2").Value = 1532.1666
$ws.Range("K2").Value = 147958.86
$ws.Range("L2").Value = 1532.1666
$ws.Range("M2").Value = -147845.86
$ws.Range("N2").Value = -1758.1666
$ws.Range("H45").Value = 1592.75
$ws.Range("I45").Value = 1014.1429
$ws.Range("K45").Value = 1014.1429
$ws.Range("M45").Value = -637.1429000000001
$ws.Range("H110").Value = 2259.647
$ws.Range("I110").Value = 1727.75
$ws.Range("J110").Value = 2423.3076
$ws.Range("K110").Value = 1727.75
$ws.Range("L110").Value = 2423.3076
$ws.Range("M110").Value = 317.25
$ws.Range("N110").Value = -6513.3076
$ws.Range("H116").Value = 80377.30499999999
$ws.Range("I116").Value = 147958.86
$ws.Range("J116").Value = 1532.1666
$ws.Range("K116").Value = 147958.86
$ws.Range("L116").Value = 1532.1666
$ws.Range("M116").Value = -145664.86
$ws.Range("N116").Value = -6120.1666
$ws.Range("H135").Value = 34298.25
$ws.Range("J135").Value = 34298.25
$ws.Range("L135").Value = 34298.25
$ws.Range("N135").Value = -44438.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 80377.30499999999
$ws.Range("I3").Value = 147958.86
$ws.Range("J3").Value = 1532.1666
$ws.Range("K3").Value = 147958.86
$ws.Range("L3").Value = 1532.1666
$ws.Range("M3").Value = -147844.86
$ws.Range("N3").Value = -1760.1666
$ws.Range("H107").Value = 1438.3125
$ws.Range("I107").Value = 1165.5834
$ws.Range("K107").Value = 1165.5834
$ws.Range("M107").Value = 754.4166

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1250
$ws.Range("I16").Value = 1250
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1250
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -963
$ws.Range("N16").ClearContents()
$ws.Range("H107").Value = 358.69565
$ws.Range("I107").Value = 187.5
$ws.Range("J107").Value = 490.3846
$ws.Range("K107").Value = 187.5
$ws.Range("L107").Value = 490.3846
$ws.Range("M107").Value = 1732.5
$ws.Range("N107").Value = -4330.3846
$ws.Range("H113").Value = 1250
$ws.Range("I113").Value = 1250
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1250
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 920
$ws.Range("N113").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1273.8096
$ws.Range("I51").Value = 2875
$ws.Range("J51").Value = 1105.2632
$ws.Range("K51").Value = 8625
$ws.Range("L51").Value = 3315.7896
$ws.Range("M51").Value = -8165
$ws.Range("N51").Value = -4235.7896
$ws.Range("H103").Value = 4121.3335
$ws.Range("I103").Value = 2000
$ws.Range("J103").Value = 4545.6
$ws.Range("K103").Value = 6000
$ws.Range("L103").Value = 13636.8
$ws.Range("M103").Value = -5121
$ws.Range("N103").Value = -15394.8
$ws.Range("H140").Value = 5818.1333
$ws.Range("I140").Value = 11604.9
$ws.Range("J140").Value = 2924.75
$ws.Range("K140").Value = 34814.7
$ws.Range("L140").Value = 8774.25
$ws.Range("M140").Value = -29634.7
$ws.Range("N140").Value = -19134.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 860.0714
$ws.Range("I107").Value = 1214.625
$ws.Range("J107").Value = 387.33334
$ws.Range("K107").Value = 1214.625
$ws.Range("L107").Value = 387.33334
$ws.Range("M107").Value = 705.375
$ws.Range("N107").Value = -4227.33334
$ws.Range("H113").Value = 1572.2667
$ws.Range("I113").Value = 666.6667
$ws.Range("J113").Value = 1798.6666
$ws.Range("K113").Value = 666.6667
$ws.Range("L113").Value = 1798.6666
$ws.Range("M113").Value = 1503.3333
$ws.Range("N113").Value = -6138.6666
$ws.Range("H126").Value = 2538.3235
$ws.Range("I126").Value = 2149.0833
$ws.Range("J126").Value = 2750.6365
$ws.Range("K126").Value = 6447.249899999999
$ws.Range("L126").Value = 8251.9095
$ws.Range("M126").Value = -3977.249899999999
$ws.Range("N126").Value = -13191.9095

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1447.4348
$ws.Range("I61").Value = 1136.3684
$ws.Range("K61").Value = 1136.3684
$ws.Range("M61").Value = -934.3684000000001
$ws.Range("H113").Value = 1447.4348
$ws.Range("I113").Value = 1136.3684
$ws.Range("K113").Value = 1136.3684
$ws.Range("M113").Value = 1033.6316

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 1826.5
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1826.5
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 1826.5
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -2380.5
$ws.Range("H107").Value = 949.6667
$ws.Range("I107").Value = 900
$ws.Range("J107").Value = 999.3333
$ws.Range("K107").Value = 2700
$ws.Range("L107").Value = 2997.9999
$ws.Range("M107").Value = -780
$ws.Range("N107").Value = -6837.9999
$ws.Range("H113").Value = 435
$ws.Range("I113").Value = 320
$ws.Range("J113").Value = 550
$ws.Range("K113").Value = 960
$ws.Range("L113").Value = 1650
$ws.Range("M113").Value = 1210
$ws.Range("N113").Value = -5990
$ws.Range("H125").Value = 28686.55
$ws.Range("J125").Value = 28686.55
$ws.Range("L125").Value = 28686.55
$ws.Range("N125").Value = -38526.55
